$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-79 from 45202 to 45203
$ws.Range("C2:C79").Value = 45203

# Add new row 80 with data
$ws.Cells.Item(80, 1).Value = "A 47231-2023"
$ws.Cells.Item(80, 2).Value = 45202
$ws.Cells.Item(80, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(80, 3).Value = 45203
$ws.Cells.Item(80, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(80, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(80, 5).Value = "TOMELILLA"
$ws.Cells.Item(80, 7).Value = 1.3
$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 13).Value = 0
$ws.Cells.Item(80, 14).Value = 0
$ws.Cells.Item(80, 15).Value = 0
$ws.Cells.Item(80, 16).Value = 0
$ws.Cells.Item(80, 17).Value = 0
$ws.Cells.Item(80, 18).Value = ""
$ws.Cells.Item(80, 18).WrapText = $true

$ws.Rows.Item(79).RowHeight = 15
